$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.353.83"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.12%  '

$ws.Range('D3').Value = "'2.090.15"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.31%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = "'252.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.94%  '

$ws.Range('E6').Value = '  +0.28%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = "'54.24"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +20.26%  '

$ws.Range('D9').Value = "'62.43"
$ws.Range('D9').Style = 'Normal'

$ws.Range('D10').Value = "'0.381"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.65%  '

$ws.Range('D11').Value = "'0.0756"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.49%  '

$ws.Range('E12').Value = '  +7.45%  '

$ws.Range('D13').Value = "'15.46"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.58%  '

$ws.Range('D14').Value = "'2.393.87"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.31%  '

$ws.Range('D15').Value = "'0.859"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.91%  '

$ws.Range('E16').Value = '  +7.51%  '

$ws.Range('D17').Value = "'2.093.18"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('D18').Value = "'37.318.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.09%  '

$ws.Range('D19').Value = "'73.23"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.18%  '

$ws.Range('E20').Value = '  +13.62%  '

$ws.Range('E21').Value = '  +4.92%  '

$ws.Range('D22').Value = "'241.20"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.08%  '

$ws.Range('E23').Value = '  +7.24%  '

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  +0.55%  '

$ws.Range('E26').Value = '  +1.94%  '

$ws.Range('E27').Value = '  +4.76%  '

$ws.Range('D28').Value = "'20.89"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.44%  '

$ws.Range('E29').Value = '  +3.71%  '

$ws.Range('E30').Value = '  +2.43%  '

$ws.Range('D31').Value = "'23.58"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.84%  '

$ws.Range('D32').Value = "'1.10"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +21.99%  '

$ws.Range('E33').Value = '  +4.44%  '

$ws.Range('E34').Value = '  +7.51%  '

$ws.Range('D35').Value = "'0.0901"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.11%  '

$ws.Range('E36').Value = '  +7.31%  '

$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = "'2.28"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.77%  '

$ws.Range('E39').Value = '  -4.13%  '

$ws.Range('D40').Value = "'1.36"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.37%  '

$ws.Range('E41').Value = '  +6.40%  '

$ws.Range('D42').Value = "'17.95"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +12.78%  '

$ws.Range('E43').Value = '  +3.54%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'99.80"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.05%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.0972"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +18.19%  '

$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = "'4.45"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +123.86%  '

$ws.Range('D47').Value = "'2.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.22%  '

$ws.Range('D48').Value = "'1.334.73"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.25%  '

$ws.Range('E49').Value = '  +3.79%  '

$ws.Range('D50').Value = "'2.35"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.44%  '

$ws.Range('D51').Value = "'6.98"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +12.61%  '
